# Update the marksheet correct/total marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: number of right answers used for scoring (B11)
$ws.Range("B11").Value = 5

# "Total" row: total score achieved (B12)
$ws.Range("B12").Value = 85

# Corr/total marks display text (E12)
$ws.Range("E12").Value = "85/140"
